# Update the "Fitness" values in column C for this run's log sheet.
# Original values were a flat 7310 for every row; replace with the
# recorded fitness values for the corresponding generation ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: start row, end row (inclusive), new value for column C
$ranges = @(
    @(2,   2,   8971),
    @(3,   4,   8947),
    @(5,   5,   8700),
    @(6,   10,  8128),
    @(11,  23,  7768),
    @(24,  42,  7754),
    @(43,  103, 7721),
    @(215, 252, 7293)
)

foreach ($r in $ranges) {
    $startRow = $r[0]
    $endRow = $r[1]
    $value = $r[2]
    $rangeAddr = "C$startRow`:C$endRow"
    $ws.Range($rangeAddr).Value = $value
}
